$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to remain plain text even when the new value
# looks numeric (e.g. "99.03", "5.00", "3.80"), matching the inlineStr cells
# in the workbook. NumberFormat="@" marks the cell as Text before the write so
# Excel does not silently coerce/round the string into a Number, then
# ClearFormats() strips the Text number-format back off again (restoring the
# original default cell style) without touching the stored string value.
function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextValue "D2" "65.290.56"
$ws.Range("E2").Value = "  -6.19%  "
Set-TextValue "D3" "3.290.54"
$ws.Range("E3").Value = "  -7.29%  "
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue "D5" "551.81"
$ws.Range("E5").Value = "  -5.81%  "
Set-TextValue "D6" "178.08"
$ws.Range("E6").Value = "  -9.22%  "
$ws.Range("E7").Value = "  +0.13%  "
Set-TextValue "D8" "0.581"
$ws.Range("E8").Value = "  -4.92%  "
Set-TextValue "D9" "3.289.39"
$ws.Range("E9").Value = "  -6.94%  "
Set-TextValue "D10" "0.181"
$ws.Range("E10").Value = "  -12.36%  "
Set-TextValue "D11" "0.576"
$ws.Range("E11").Value = "  -8.75%  "
Set-TextValue "D12" "46.54"
$ws.Range("E12").Value = "  -11.94%  "
Set-TextValue "D13" "0.0000259"
$ws.Range("E13").Value = "  -10.20%  "
Set-TextValue "D14" "3.833.19"
$ws.Range("E14").Value = "  -6.89%  "
Set-TextValue "D15" "8.41"
$ws.Range("E15").Value = "  -9.02%  "
Set-TextValue "D16" "589.40"
$ws.Range("E16").Value = "  -11.54%  "
Set-TextValue "D17" "65.364.22"
$ws.Range("E17").Value = "  -6.27%  "
Set-TextValue "D18" "17.77"
$ws.Range("E18").Value = "  -3.68%  "
$ws.Range("E19").Value = "  -4.26%  "
Set-TextValue "D20" "3.304.93"
$ws.Range("E20").Value = "  -7.19%  "
Set-TextValue "D21" "11.22"
$ws.Range("E21").Value = "  -10.28%  "
Set-TextValue "D22" "0.887"
$ws.Range("E22").Value = "  -8.13%  "
Set-TextValue "D23" "16.63"
$ws.Range("E23").Value = "  -7.79%  "
Set-TextValue "D24" "5.00"
$ws.Range("E24").Value = "  -6.72%  "
Set-TextValue "D25" "99.03"
$ws.Range("E25").Value = "  -5.96%  "
Set-TextValue "D26" "3.93"
$ws.Range("E26").Value = "  -10.40%  "
Set-TextValue "D27" "5.97"
$ws.Range("E27").Value = "  -0.51%  "
Set-TextValue "D28" "2.62"
$ws.Range("E28").Value = "  -10.45%  "
Set-TextValue "D29" "9.09"
$ws.Range("E29").Value = "  -10.46%  "
Set-TextValue "D30" "8.51"
$ws.Range("E30").Value = "  -11.26%  "
Set-TextValue "D31" "30.06"
$ws.Range("E31").Value = "  -9.96%  "
Set-TextValue "D32" "3.80"
$ws.Range("E32").Value = "  -13.00%  "
Set-TextValue "D33" "6.14"
$ws.Range("E33").Value = "  -9.37%  "
Set-TextValue "D34" "10.86"
$ws.Range("E34").Value = "  -7.78%  "
$ws.Range("E35").Value = "  -7.98%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D36" "57.33"
$ws.Range("E36").Value = "  -7.56%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D37" "3.723.43"
$ws.Range("E37").Value = "  -1.89%  "
$ws.Range("E38").Value = "  -0.24%  "
Set-TextValue "D39" "515.39"
$ws.Range("E39").Value = "  +3.27%  "
Set-TextValue "D40" "3.43"
$ws.Range("E40").Value = "  -8.84%  "
Set-TextValue "D41" "0.0₃0698"
$ws.Range("E41").Value = "  -13.71%  "
Set-TextValue "D42" "2.61"
$ws.Range("E42").Value = "  -9.78%  "
Set-TextValue "D43" "0.123"
$ws.Range("E43").Value = "  -8.60%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D44" "0.333"
$ws.Range("E44").Value = "  -10.20%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D45" "31.33"
$ws.Range("E45").Value = "  -9.51%  "
Set-TextValue "D46" "3.28"
$ws.Range("E46").Value = "  -3.08%  "
Set-TextValue "D47" "0.0406"
$ws.Range("E47").Value = "  -10.07%  "
Set-TextValue "D48" "2.96"
$ws.Range("E48").Value = "  +9.86%  "
Set-TextValue "D49" "0.127"
$ws.Range("E49").Value = "  -6.90%  "
$ws.Range("E50").Value = "  -11.23%  "
$ws.Range("E51").Value = "  -0.06%  "
